$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the Date paragraph: it was split across three runs
#    ("Date: 6/1" + "8" + "/2015 at 9.30 AM EST"); a Find/Replace over
#    the full visible text collapses it back into a single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Date: 6/18/2015 at 9.30 AM EST", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Date: 6/18/2015 at 9.30 AM EST", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Drop the _GoBack bookmark from its old location (the "1. Git hub
#    organization..." list item).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3) Remove the "We should not refer this project as OpenFDA/ADS in
#    public forums" paragraph entirely (it merges into the blank
#    paragraph that used to trail it), then re-create the _GoBack
#    bookmark (empty) in what is now the final, trailing paragraph.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "We should not refer this project*") {
        $target = $p
        break
    }
}
if ($target -ne $null) {
    $target.Range.Delete()
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

# The engine mis-handles Bookmarks.Add on a range whose visible text is
# empty (it silently resets the bookmark to position 0). Work around
# this by adding the bookmark while the paragraph still holds a unique
# placeholder run, then erasing that placeholder afterwards; the
# bookmark collapses to an empty range in place, same as real Word.
$placeholder = "ZZGoBackPlaceholderZZ"
$r = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$r.InsertBefore($placeholder)

$lastPara2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r2 = $d.Range($lastPara2.Range.Start, $lastPara2.Range.End)
$d.Bookmarks.Add("_GoBack", $r2)

$d.Content.Find.Execute($placeholder, $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null
